$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "amount in words" line
$ws.Range("T7").Value = "(Thành tiền bằng chữ: Không đồng)"

# Update main title (was a rich-text run, now plain text)
$ws.Range("A2").Value = "BẢNG THÙ LAO TIN, PS TRONG THÔNG TIN NGÀY MỚI"

# Total amount cell now has an explicit 0
$ws.Range("S6").Value = 0

# Update active selection
$ws.Range("E10").Select()
